# Weekly update: insert a new observation row for Betarraga (Macroferia Regional
# de Talca) at row 260, pushing the existing historical rows (260-363) down by
# one row so the new week's record sits at the top of that block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 260; existing rows 260..363 shift to 261..364
$ws.Rows.Item(260).Insert()

# Populate the newly inserted row with this week's data
$ws.Cells.Item(260, 1).Value  = 5
$ws.Cells.Item(260, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(260, 3).Value  = "Maule"
$ws.Cells.Item(260, 4).Value  = 44784
$ws.Cells.Item(260, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(260, 5).Value  = 7
$ws.Cells.Item(260, 6).Value  = 100114014
$ws.Cells.Item(260, 7).Value  = "Betarraga"
$ws.Cells.Item(260, 8).Value  = "Sin especificar"
$ws.Cells.Item(260, 9).Value  = "Primera"
$ws.Cells.Item(260, 10).Value = 3000
$ws.Cells.Item(260, 11).Value = 750
$ws.Cells.Item(260, 12).Value = 750
$ws.Cells.Item(260, 13).Value = 750
$ws.Cells.Item(260, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(260, 15).Value = "Región del Maule"
$ws.Cells.Item(260, 16).Value = 150
$ws.Cells.Item(260, 17).Value = 5
$ws.Cells.Item(260, 18).Value = "Hortaliza"
